$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting the existing row 19 (and below) down to row 20.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new match data.
$ws.Range("A19").Value = "tYtqElCi"
$ws.Range("B19").Value = "17/11/2024"
$ws.Range("C19").Value = "18:30"
$ws.Range("D19").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E19").Value = "Penarol"
$ws.Range("F19").Value = "Defensor Sp."
$ws.Range("G19").Value = 1.44
$ws.Range("H19").Value = 4.5
$ws.Range("I19").Value = 7
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 2.2
$ws.Range("L19").Value = 8.5
$ws.Range("M19").Value = 1.06
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 1.4
$ws.Range("P19").Value = 2.75
$ws.Range("Q19").Value = 2.2
$ws.Range("R19").Value = 1.65
$ws.Range("S19").Value = 1.44
$ws.Range("T19").Value = 2.63
$ws.Range("U19").Value = 2.5
$ws.Range("V19").Value = 1.5
$ws.Range("W19").Value = 5
$ws.Range("X19").Value = 5.5
$ws.Range("Y19").Value = 9.5
$ws.Range("Z19").Value = 8.5
$ws.Range("AA19").Value = 15
$ws.Range("AB19").Value = 41
$ws.Range("AC19").Value = 8.5
$ws.Range("AD19").Value = 9
$ws.Range("AE19").Value = 29
$ws.Range("AF19").Value = 101
$ws.Range("AG19").Value = 13
$ws.Range("AH19").Value = 34
$ws.Range("AI19").Value = 23
$ws.Range("AJ19").Value = 101
$ws.Range("AK19").Value = 67
$ws.Range("AL19").Value = 67
$ws.Range("AM19").Value = 201
$ws.Range("AN19").Value = 3.1
$ws.Range("AO19").Value = 7
$ws.Range("AP19").Value = 23
$ws.Range("AQ19").Value = 23
$ws.Range("AR19").Value = 51
$ws.Range("AS19").Value = 251
$ws.Range("AT19").Value = 2.63
$ws.Range("AU19").Value = 11
$ws.Range("AV19").Value = 81
$ws.Range("AW19").Value = 9
$ws.Range("AX19").Value = 41
$ws.Range("AY19").Value = 51
$ws.Range("AZ19").Value = 251
$ws.Range("BA19").Value = 301
$ws.Range("BB19").Value = 501
$ws.Range("BC19").Value = 51
$ws.Range("BD19").Value = 51

# The row that used to be row 19 is now row 20; update its one changed value.
$ws.Range("N20").Value = 6.2

# Apply the odds updates for the other existing matches (rows 4, 5, 8, 10, 11, 14, 15, 16, 17).
$ws.Range("G4").Value = 2.3
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 3.2
$ws.Range("J4").Value = 3.1
$ws.Range("N4").Value = 7.5
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("Y4").Value = 10
$ws.Range("AG4").Value = 8
$ws.Range("AN4").Value = 4.33
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 81
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.38
$ws.Range("AY4").Value = 29
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("BD8").Value = 126
$ws.Range("G10").Value = 3.7
$ws.Range("I10").Value = 1.91
$ws.Range("J10").Value = 4.75
$ws.Range("L10").Value = 2.63
$ws.Range("X10").Value = 19
$ws.Range("AC10").Value = 8.5
$ws.Range("AG10").Value = 6
$ws.Range("AH10").Value = 8
$ws.Range("AJ10").Value = 15
$ws.Range("AN10").Value = 6
$ws.Range("AR10").Value = 126
$ws.Range("AS10").Value = 301
$ws.Range("AU10").Value = 9
$ws.Range("AW10").Value = 3.75
$ws.Range("G11").Value = 4.1
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 1.9
$ws.Range("J11").Value = 4.75
$ws.Range("L11").Value = 2.6
$ws.Range("X11").Value = 21
$ws.Range("Y11").Value = 15
$ws.Range("AH11").Value = 8.5
$ws.Range("AJ11").Value = 15
$ws.Range("AN11").Value = 6
$ws.Range("AO11").Value = 23
$ws.Range("AW11").Value = 3.75
$ws.Range("BA11").Value = 51
$ws.Range("G14").Value = 1.73
$ws.Range("H14").Value = 3.8
$ws.Range("I14").Value = 4.5
$ws.Range("J14").Value = 2.25
$ws.Range("AH14").Value = 26
$ws.Range("AJ14").Value = 51
$ws.Range("AK14").Value = 34
$ws.Range("AL14").Value = 34
$ws.Range("AQ14").Value = 26
$ws.Range("AW14").Value = 6.5
$ws.Range("M15").Value = 1.11
$ws.Range("N15").Value = 6.5
$ws.Range("G16").Value = 2.05
$ws.Range("I16").Value = 4.2
$ws.Range("I17").Value = 2.5
$ws.Range("L17").Value = 3.25
$ws.Range("AA17").Value = 23
$ws.Range("AD17").Value = 6
$ws.Range("AM17").Value = 301
$ws.Range("BA17").Value = 81
